$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 53: date (stored as plain text, matching the rest of the column)
# and the numeric profit value.
$dateCell = $ws.Cells.Item(53, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/09/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item(53, 2).Value = 14480.62
